$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet 1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1885
$ws1.Range("F3").Value = 1514
$ws1.Range("F4").Value = 882
$ws1.Range("F5").Value = 767
$ws1.Range("F6").Value = 13328
$ws1.Range("F7").Value = 13190
$ws1.Range("F8").Value = 1019
$ws1.Range("F9").Value = 772
$ws1.Range("F10").Value = 22
$ws1.Range("F11").Value = 557
$ws1.Range("F13").Value = 675
$ws1.Range("F14").Value = 2095
$ws1.Range("F19").Value = 398
$ws1.Range("F20").Value = 254
$ws1.Range("F21").Value = 288
$ws1.Range("F22").Value = 421
$ws1.Range("F23").Value = 757
$ws1.Range("F24").Value = 16

# Sheet "演出" (sheet 2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 61
$ws2.Range("F7").Value = 120
$ws2.Range("F9").Value = 30

# Sheet "本地生活" (sheet 3)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 193
$ws3.Range("F3").Value = 41

# Sheet "全部类型" (sheet 4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 193
$ws4.Range("F3").Value = 1885
$ws4.Range("F4").Value = 1514
$ws4.Range("F5").Value = 882
$ws4.Range("F7").Value = 767
$ws4.Range("F8").Value = 13328
$ws4.Range("F9").Value = 13190
$ws4.Range("F10").Value = 1019
$ws4.Range("F11").Value = 772
$ws4.Range("F12").Value = 22
$ws4.Range("F13").Value = 557
$ws4.Range("F15").Value = 675
$ws4.Range("F18").Value = 2095
$ws4.Range("F24").Value = 61
$ws4.Range("F25").Value = 41
$ws4.Range("F26").Value = 398
$ws4.Range("F27").Value = 254
$ws4.Range("F28").Value = 288
$ws4.Range("F29").Value = 421
$ws4.Range("F30").Value = 757
$ws4.Range("F31").Value = 120
$ws4.Range("F33").Value = 16
$ws4.Range("F34").Value = 30
